# Auto-generated Word COM-interop script.
# Adds the expanded citation-list paragraph + full bibliography entries
# described by the commit "Added to bibliography, incorporated Matt's
# changes, formatted tables for Word".

$d = $word.ActiveDocument

# Helper: replace the *entire* text of a paragraph (all runs) safely.
# (Assigning directly to $para.Range.Text only clobbers the first run
#  when the paragraph has more than one run, so we build an explicit
#  sub-range that excludes the trailing paragraph mark.)
function Set-ParaText($para, $text) {
    $r = $para.Range
    $ir = $r.Duplicate
    $ir.SetRange($r.Start, $r.End - 1)
    $ir.Text = $text
}

# Helper: append another run of text onto the end of a paragraph,
# right before its paragraph mark.
function Append-ToPara($para, $text) {
    $r = $para.Range
    $ir = $r.Duplicate
    $ir.SetRange($r.End - 1, $r.End - 1)
    $ir.InsertAfter($text)
}

# --- Step 1: insert the new intro paragraph (expanded citation list) ---
# --- right before the "References" heading.                          ---
$headingPara = $d.Paragraphs.Item(2)
$headingPara.Range.InsertParagraphBefore()
$newFirst = $d.Paragraphs.Item(2)
$newFirst.Style = "FirstParagraph"
Set-ParaText $newFirst 'Blah blah'
Append-ToPara $newFirst ' '
Append-ToPara $newFirst '(Balesdent and Balabane, 1996; Beniston et al., 2014; Blackmer, 1997; Buyanovsky et al., 1987; Cotrufo et al., 2015; David et al., 2009; Davidson and Ackerman, 1993; Dupont et al., 2014; Gill and Burke, 2002; Gill et al., 1999; Gregory et al., 2016; Guo and Gifford, 2002; Guzman and Al-Kaisi, 2010; Heggenstaller et al., 2009; Huggins et al., 1998; Jarchow and Liebman, 2013; Jobbágy and Jackson, 2000; Kong and Six, 2010; Liang and Balser, 2008; McGranahan et al., 2014; Omonode and Vyn, 2006; O’BRIEN et al., 2010; Pinheiro et al., 2013; Rasse et al., 2005; Rumpel and Kögel-Knabner, 2011; Silver and Miya, 2001; Six et al., 2002; Tufekcioglu et al., 2003; Van Es et al., 2007; abendroth2011; Wiles et al., 1996)'
Append-ToPara $newFirst '.'

# --- Step 2: turn the old intro paragraph into the first batch of new ---
# --- Bibliography entries (Balesdent .. Huggins), inserted right      ---
# --- after "References" and before the pre-existing "Jarchow" entry.  ---
$oldFirst = $d.Paragraphs.Item(4)
$oldFirst.Style = "Bibliography"
Set-ParaText $oldFirst 'Balesdent, J. and Balabane, M.: Major contribution of roots to soil carbon storage inferred from maize cultivated soils, Soil Biology and Biochemistry, 28(9), 1261–1263, 1996.'
$prev = $oldFirst
$prev.Range.InsertParagraphAfter()
$prevIdx = $prev.Index
$prev = $d.Paragraphs.Item($prevIdx + 1)
Set-ParaText $prev 'Beniston, J. W., DuPont, S. T., Glover, J. D., Lal, R. and Dungait, J. A.: Soil organic carbon dynamics 75 years after land-use change in perennial grassland and annual wheat agricultural systems, Biogeochemistry, 120(1-3), 37–49, 2014.'
$prev.Range.InsertParagraphAfter()
$prevIdx = $prev.Index
$prev = $d.Paragraphs.Item($prevIdx + 1)
Set-ParaText $prev 'Blackmer, A.: Nitrogen fertilizer recommendations for corn in iowa, 1997.'
$prev.Range.InsertParagraphAfter()
$prevIdx = $prev.Index
$prev = $d.Paragraphs.Item($prevIdx + 1)
Set-ParaText $prev 'Buyanovsky, G., Kucera, C. and Wagner, G.: Comparative analyses of carbon dynamics in native and cultivated ecosystems, Ecology, 68(6), 2023–2031, 1987.'
$prev.Range.InsertParagraphAfter()
$prevIdx = $prev.Index
$prev = $d.Paragraphs.Item($prevIdx + 1)
Set-ParaText $prev 'Cotrufo, M. F., Soong, J. L., Horton, A. J., Campbell, E. E., Haddix, M. L., Wall, D. H. and Parton, W. J.: Formation of soil organic matter via biochemical and physical pathways of litter mass loss, Nature Geoscience, 2015.'
$prev.Range.InsertParagraphAfter()
$prevIdx = $prev.Index
$prev = $d.Paragraphs.Item($prevIdx + 1)
Set-ParaText $prev 'David, M. B., McIsaac, G. F., Darmody, R. G. and Omonode, R. A.: Long-term changes in mollisol organic carbon and nitrogen, Journal of Environmental Quality, 38(1), 200–211, 2009.'
$prev.Range.InsertParagraphAfter()
$prevIdx = $prev.Index
$prev = $d.Paragraphs.Item($prevIdx + 1)
Set-ParaText $prev 'Davidson, E. A. and Ackerman, I. L.: Changes in soil carbon inventories following cultivation of previously untilled soils, Biogeochemistry, 20(3), 161–193, 1993.'
$prev.Range.InsertParagraphAfter()
$prevIdx = $prev.Index
$prev = $d.Paragraphs.Item($prevIdx + 1)
Set-ParaText $prev 'Dupont, S. T., Beniston, J., Glover, J., Hodson, A., Culman, S., Lal, R. and Ferris, H.: Root traits and soil properties in harvested perennial grassland, annual wheat, and never-tilled annual wheat, Plant and soil, 381(1-2), 405–420, 2014.'
$prev.Range.InsertParagraphAfter()
$prevIdx = $prev.Index
$prev = $d.Paragraphs.Item($prevIdx + 1)
Set-ParaText $prev 'Gill, R., Burke, I. C., Milchunas, D. G. and Lauenroth, W. K.: Relationship between root biomass and soil organic matter pools in the shortgrass steppe of eastern colorado, Ecosystems, 2(3), 226–236, 1999.'
$prev.Range.InsertParagraphAfter()
$prevIdx = $prev.Index
$prev = $d.Paragraphs.Item($prevIdx + 1)
Set-ParaText $prev 'Gill, R. A. and Burke, I. C.: Influence of soil depth on the decomposition of bouteloua gracilis roots in the shortgrass steppe, Plant and Soil, 241(2), 233–242, 2002.'
$prev.Range.InsertParagraphAfter()
$prevIdx = $prev.Index
$prev = $d.Paragraphs.Item($prevIdx + 1)
Set-ParaText $prev 'Gregory, A., Dungait, J., Watts, C., Bol, R., Dixon, E., White, R. and Whitmore, A.: Long-term management changes topsoil and subsoil organic carbon and nitrogen dynamics in a temperate agricultural system, European journal of soil science, 67(4), 421–430, 2016.'
$prev.Range.InsertParagraphAfter()
$prevIdx = $prev.Index
$prev = $d.Paragraphs.Item($prevIdx + 1)
Set-ParaText $prev 'Guo, L. B. and Gifford, R.: Soil carbon stocks and land use change: A meta analysis, Global change biology, 8(4), 345–360, 2002.'
$prev.Range.InsertParagraphAfter()
$prevIdx = $prev.Index
$prev = $d.Paragraphs.Item($prevIdx + 1)
Set-ParaText $prev 'Guzman, J. G. and Al-Kaisi, M. M.: Soil carbon dynamics and carbon budget of newly reconstructed tall-grass prairies in south central iowa, Journal of environmental quality, 39(1), 136–146, 2010.'
$prev.Range.InsertParagraphAfter()
$prevIdx = $prev.Index
$prev = $d.Paragraphs.Item($prevIdx + 1)
Set-ParaText $prev 'Heggenstaller, A. H., Moore, K. J., Liebman, M. and Anex, R. P.: Nitrogen influences biomass and nutrient partitioning by perennial, warm-season grasses, Agronomy Journal, 101(6), 1363–1371, 2009.'
$prev.Range.InsertParagraphAfter()
$prevIdx = $prev.Index
$prev = $d.Paragraphs.Item($prevIdx + 1)
Set-ParaText $prev 'Huggins, D., Buyanovsky, G., Wagner, G., Brown, J., Darmody, R., Peck, T., Lesoing, G., Vanotti, M. and Bundy, L.: Soil organic c in the tallgrass prairie-derived region of the corn belt: Effects of long-term crop management, Soil and Tillage Research, 47(3), 219–234, 1998.'

# --- Step 3: insert more Bibliography entries (Jobbagy .. Tufekcioglu) ---
# --- right after the pre-existing "Jarchow" entry and before the       ---
# --- pre-existing "Van Es" entry (both of which are unchanged).        ---
$jarchowPara = $d.Paragraphs.Item($prev.Index + 1)
$prev2 = $jarchowPara
$prev2.Range.InsertParagraphAfter()
$prev2Idx = $prev2.Index
$prev2 = $d.Paragraphs.Item($prev2Idx + 1)
Set-ParaText $prev2 'Jobbágy, E. G. and Jackson, R. B.: The vertical distribution of soil organic carbon and its relation to climate and vegetation, Ecological applications, 10(2), 423–436, 2000.'
$prev2.Range.InsertParagraphAfter()
$prev2Idx = $prev2.Index
$prev2 = $d.Paragraphs.Item($prev2Idx + 1)
Set-ParaText $prev2 'Kong, A. Y. and Six, J.: Tracing root vs. residue carbon into soils from conventional and alternative cropping systems, Soil Science Society of America Journal, 74(4), 1201–1210, 2010.'
$prev2.Range.InsertParagraphAfter()
$prev2Idx = $prev2.Index
$prev2 = $d.Paragraphs.Item($prev2Idx + 1)
Set-ParaText $prev2 'Liang, C. and Balser, T. C.: Preferential sequestration of microbial carbon in subsoils of a glacial-landscape toposequence, dane county, wI, uSA, Geoderma, 148(1), 113–119, 2008.'
$prev2.Range.InsertParagraphAfter()
$prev2Idx = $prev2.Index
$prev2 = $d.Paragraphs.Item($prev2Idx + 1)
Set-ParaText $prev2 'McGranahan, D. A., Daigh, A. L., Veenstra, J. J., Engle, D. M., Miller, J. R. and Debinski, D. M.: Connecting soil organic carbon and root biomass with land-use and vegetation in temperate grassland, The Scientific World Journal, 2014, 2014.'
$prev2.Range.InsertParagraphAfter()
$prev2Idx = $prev2.Index
$prev2 = $d.Paragraphs.Item($prev2Idx + 1)
Set-ParaText $prev2 'Omonode, R. A. and Vyn, T. J.: Vertical distribution of soil organic carbon and nitrogen under warm-season native grasses relative to croplands in west-central indiana, uSA, Agriculture, Ecosystems & Environment, 117(2), 159–170, 2006.'
$prev2.Range.InsertParagraphAfter()
$prev2Idx = $prev2.Index
$prev2 = $d.Paragraphs.Item($prev2Idx + 1)
Set-ParaText $prev2 'O’BRIEN, S. L., Jastrow, J. D., Grimley, D. A. and GONZALEZ-MELER, M. A.: Moisture and vegetation controls on decadal-scale accrual of soil organic carbon and total nitrogen in restored grasslands, Global Change Biology, 16(9), 2573–2588, 2010.'
$prev2.Range.InsertParagraphAfter()
$prev2Idx = $prev2.Index
$prev2 = $d.Paragraphs.Item($prev2Idx + 1)
Set-ParaText $prev2 'Pinheiro, J., Bates, D., DebRoy, S., Sarkar, D. and Team, R. C.: Nlme: Linear and nonlinear mixed effects models. r package version 3.1-113, available at h ttp, CRAN. R-project. org/package= nlme, 2013.'
$prev2.Range.InsertParagraphAfter()
$prev2Idx = $prev2.Index
$prev2 = $d.Paragraphs.Item($prev2Idx + 1)
Set-ParaText $prev2 'Rasse, D. P., Rumpel, C. and Dignac, M.-F.: Is soil carbon mostly root carbon? Mechanisms for a specific stabilisation, Plant and soil, 269(1-2), 341–356, 2005.'
$prev2.Range.InsertParagraphAfter()
$prev2Idx = $prev2.Index
$prev2 = $d.Paragraphs.Item($prev2Idx + 1)
Set-ParaText $prev2 'Rumpel, C. and Kögel-Knabner, I.: Deep soil organic matter-a key but poorly understood component of terrestrial c cycle, Plant and Soil, 338(1-2), 143–158, 2011.'
$prev2.Range.InsertParagraphAfter()
$prev2Idx = $prev2.Index
$prev2 = $d.Paragraphs.Item($prev2Idx + 1)
Set-ParaText $prev2 'Silver, W. L. and Miya, R. K.: Global patterns in root decomposition: Comparisons of climate and litter quality effects, Oecologia, 129(3), 407–419, 2001.'
$prev2.Range.InsertParagraphAfter()
$prev2Idx = $prev2.Index
$prev2 = $d.Paragraphs.Item($prev2Idx + 1)
Set-ParaText $prev2 'Six, J., Conant, R., Paul, E. A. and Paustian, K.: Stabilization mechanisms of soil organic matter: Implications for c-saturation of soils, Plant and soil, 241(2), 155–176, 2002.'
$prev2.Range.InsertParagraphAfter()
$prev2Idx = $prev2.Index
$prev2 = $d.Paragraphs.Item($prev2Idx + 1)
Set-ParaText $prev2 'Tufekcioglu, A., Raich, J., Isenhart, T. and Schultz, R.: Biomass, carbon and nitrogen dynamics of multi-species riparian buffers within an agricultural watershed in iowa, uSA, Agroforestry Systems, 57(3), 187–198, 2003.'

# --- Step 4: insert the final Bibliography entry (Wiles) right after  ---
# --- the pre-existing "Van Es" entry.                                  ---
$vanEsPara = $d.Paragraphs.Item($prev2.Index + 1)
$vanEsPara.Range.InsertParagraphAfter()
$wilesPara = $d.Paragraphs.Item($vanEsPara.Index + 1)
Set-ParaText $wilesPara 'Wiles, L. J., Barlin, D. H., Schweizer, E. E., Duke, H. R. and Whitt, D. E.: A new soil sampler and elutriator for collecting and extracting weed seeds from soil, Weed Technology, 35–41, 1996.'

